# "Move feature eng to page1 and add list"
#
# 1. Data-Check sheet: "Missing value imputation" -> "Missing value validation"
# 2. Add new sheet "Feature Engineering" (after Data-Check)
# 3. Add new sheet "Data-Clean" (after Feature Engineering)
# 4. Populate both new sheets, and set the active tab to Data-Clean.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Data-Check: rewrite the "Missing value imputation" row as
#    "Missing value validation"
# ---------------------------------------------------------------------------
$wsCheck = $wb.Worksheets.Item(1)
$wsCheck.Range("A3").Value = "Missing value validation"

# ---------------------------------------------------------------------------
# 2) Add "Feature Engineering" sheet right after Data-Check
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFeat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$wsFeat.Name = "Feature Engineering"

# Column B (time domain stats) first - matches original authoring order
$wsFeat.Range("B1").Value = "mean "
$wsFeat.Range("B2").Value = "median"
$wsFeat.Range("B3").Value = "standard deviation"
$wsFeat.Range("B4").Value = "variance"
$wsFeat.Range("B5").Value = "skewness"
$wsFeat.Range("B6").Value = "kurtosis"
$wsFeat.Range("B7").Value = "square root"
$wsFeat.Range("B8").Value = "peak-to-peak range"
$wsFeat.Range("B9").Value = "root mean square"
$wsFeat.Range("B10").Value = "sum of squares"

# Section headers (column A)
$wsFeat.Range("A11").Value = "Frequency Domain Feature"
$wsFeat.Range("A1").Value = "Time Domain Feature"

# Frequency domain feature rows (column B)
$wsFeat.Range("B11").Value = "FFT (Fast Fourier Transform)"
$wsFeat.Range("B12").Value = "FFT magnitude"
$wsFeat.Range("B13").Value = "FFT frequency"
$wsFeat.Range("B14").Value = "Power spectrum"
$wsFeat.Range("B15").Value = "FFT mean"
$wsFeat.Range("B16").Value = "FFT standard deviation"
$wsFeat.Range("B17").Value = "FFT maximum"
$wsFeat.Range("B18").Value = "FFT frequency of maximum amplitude"
$wsFeat.Range("B19").Value = "Spectral centroid"
$wsFeat.Range("B20").Value = "Spectral bandwidth"
$wsFeat.Range("B21").Value = "Amplitude envelope"
$wsFeat.Range("B22").Value = "Phase envelope"
$wsFeat.Range("B23").Value = "Log power spectrum"
$wsFeat.Range("B24").Value = "Cepstrum"
$wsFeat.Range("B25").Value = "Cepstrum mean"
$wsFeat.Range("B26").Value = "Cepstrum standard deviation"
$wsFeat.Range("B27").Value = "Cepstrum maximum"

# Time-frequency section header + rows
$wsFeat.Range("A28").Value = "Time Frequency Feature"
$wsFeat.Range("B28").Value = "STFT mean"
$wsFeat.Range("B29").Value = "STFT standard deviation"
$wsFeat.Range("B30").Value = "STFT maximum"
$wsFeat.Range("B31").Value = "Wavelet mean"
$wsFeat.Range("B32").Value = "Wavelet standard deviation"
$wsFeat.Range("B33").Value = "Wavelet maximum"

[void]$wsFeat.Columns.Item(1).AutoFit()
[void]$wsFeat.Range("A1:B10").Select()

# ---------------------------------------------------------------------------
# 3) Add "Data-Clean" sheet right after Feature Engineering
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsClean = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$wsClean.Name = "Data-Clean"

$wsClean.Range("A2").Value = "Missing Value Imputation"
$wsClean.Range("A3").Value = "Normalization"
$wsClean.Range("A4").Value = "Label Encoding"
$wsClean.Range("A1").Value = "Column Type Regulation"
$wsClean.Range("A5").Value = "Outlier Removal"

[void]$wsClean.Range("A1:A5").Select()

# ---------------------------------------------------------------------------
# 4) Selections / active tab
# ---------------------------------------------------------------------------
[void]$wsCheck.Range("B20").Select()

[void]$wsClean.Activate()
